# Updated cryptos list with GitHub Actions.
# For cells whose new text looks like a plain number (e.g. "1.001"), we force
# the cell to stay text (NumberFormat "@") and then reset the style back to
# "Normal" afterwards so no extra style index is left attached to the cell,
# matching the original inline-string, unstyled cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.127.12"
$ws.Range("D3").Value = "1.801.28"
$ws.Range("E3").Value = "  -0.25%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.37%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "311.22"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("E6").Value = "  -0.39%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5098"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -2.92%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3874"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.40%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07746"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.23%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.097"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.02%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "40.89"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.15%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "6.331"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("D15").Value = "1.797.09"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("E16").Value = "  -1.04%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "92.03"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.70%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001072"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.62%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06569"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("E20").Value = "  -0.40%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.21"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.28%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.952"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "28.143.46"
$ws.Range("E23").Value = "  -0.57%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.05"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.78%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.246"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.99%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.429"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "2.005.93"
$ws.Range("E28").Value = "  -0.34%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "20.18"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.80%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "127.00"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.98%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.1090"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.40%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.045"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.55%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.648"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.57%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.527"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.00%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.07028"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.70%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "9.048"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.32%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02339"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.08%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2157"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.80%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.008"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.75%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "11.49"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -5.29%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.6098"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("E42").Value = "  -0.36%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.154"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.294"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -6.34%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5906"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.44%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.720"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.23%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "125.14"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.61%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.197"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.34%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.897"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("E51").Value = "  -1.24%  "
